$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L ("spectraltype_esphs_dr2"); everything from the
# old L onward (dr2_source_id, dr3_source_id, ...) shifts one column right.
$ws.Columns("L:L").Insert()

# New header cells.
$ws.Range("L1").Value = "spectraltype_esphs_dr2"
$ws.Range("AA1").Value = "spectraltype_esphs_dr3"

# AA1 is a brand-new cell past the old used range, so it doesn't inherit the
# bold/centered/bordered header look the Insert() above carried along for
# L1 automatically. Match it to the rest of row 1 explicitly.
$ws.Range("AA1").Font.Bold = $true
$ws.Range("AA1").HorizontalAlignment = -4108
$ws.Range("AA1").VerticalAlignment = -4160
$ws.Range("AA1").Borders.LineStyle = 1

# Gaia ESP-HS spectral-type values for the DR2 cross-match (only the rows
# that had a match in the diff get a value; others stay blank).
$ws.Range("L2").Value = "K"
$ws.Range("L3").Value = "K"
$ws.Range("L5").Value = "K"
$ws.Range("L6").Value = "K"
$ws.Range("L7").Value = "K"
$ws.Range("L8").Value = "G"
$ws.Range("L9").Value = "K"
$ws.Range("L11").Value = "K"
$ws.Range("L14").Value = "K"

# Column widths: new column L (24 chars wide) and the appended column AA
# (24 chars wide). ColumnWidth is in "characters"; the engine adds a fixed
# 5/6 character padding when it stores the OOXML <col width>, so back that
# off here to land on the exact width the diff expects.
$ws.Columns("L:L").ColumnWidth = 23.16666666666667
$ws.Columns("AA:AA").ColumnWidth = 23.16666666666667
